$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C
$ws.Range("C1").Value = "shape"

# Existing rows 2-4 get a new "circle" value in column C, and two more
# rows (5-7) are appended repeating the Red/Green/Blue + Id pattern.
$ws.Range("C2").Value = "circle"
$ws.Range("C3").Value = "circle"
$ws.Range("C4").Value = "circle"

$ws.Range("A5").Value = "Red"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "circle"

$ws.Range("A6").Value = "Green"
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = "circle"

$ws.Range("A7").Value = "Blue"
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = "circle"
